# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
#
# This workbook is a single "Estado de Cuenta" (account statement) sheet
# that lists, for a given employer, one row per worker/overdue-period
# ("Periodo Mora"). The edit adds a new overdue period (2509) for the
# existing worker, duplicating the data row (16) that already exists for
# period 2508, and updates the summary figures ("VALOR MORA" total and
# "Cant. Periodos") to reflect the additional period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new detail row (17) right below the existing data row (16),
#     duplicating its formatting, then fill in its values for the new
#     period "2509". This pushes the trailing "firma" rows down by one. ---
$ws.Rows("17:17").Insert()

$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B17").Value2 = $ws.Range("B16").Value2
$ws.Range("C17").Value2 = $ws.Range("C16").Value2
$ws.Range("D17").Value2 = $ws.Range("D16").Value2
$ws.Range("E17").Value2 = "2509"
$ws.Range("F17").Value2 = $ws.Range("F16").Value2
$ws.Range("G17").Value2 = $ws.Range("G16").Value2

# --- Update the summary block above the table: the total overdue amount
#     now covers two periods instead of one, and the period counter
#     increases from 1 to 2. ---
$ws.Range("E11").Value2 = 113880
$ws.Range("F13").Value2 = 2
